# Update "Carbon dioxide, non-fossil" rows (13-17) amounts from 0 to 2.69E-11
$ws = $excel.ActiveWorkbook.ActiveSheet

$val = 0.0000000000269

$ws.Range("C13").Value = $val
$ws.Range("C14").Value = $val
$ws.Range("C15").Value = $val
$ws.Range("C16").Value = $val
$ws.Range("C17").Value = $val

# Add two new flow rows to account for NETs (negative emissions technologies)
$ws.Range("A44").Value = "Carbon dioxide, in air"
$ws.Range("B44").Value = "natural resource::in air"
$ws.Range("C44").Value = -0.0000000000269

$ws.Range("A45").Value = "Carbon dioxide, non-fossil, resource correction"
$ws.Range("B45").Value = "natural resource::in air"
$ws.Range("C45").Value = -0.0000000000269

# Re-apply AutoFilter over the full data range (adds the hidden _FilterDatabase defined name)
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $ws.Range("A1:C45"))
$filterName.Visible = $false

# Move the active selection (matches the saved cursor position in the diff)
$ws.Range("B8").Select()
